# Scheduled runner update: refresh market-board derived profit figures
# (currentAveragePrice / NQ / HQ price & profit columns) across all Job sheets.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 1932.6666
$ws.Range("I2").Value = 924.375
$ws.Range("K2").Value = 924.375
$ws.Range("M2").Value = -811.375

$ws.Range("H40").Value = 1425
$ws.Range("J40").Value = 1425
$ws.Range("L40").Value = 1425
$ws.Range("N40").Value = -1775

$ws.Range("H58").Value = 4615.857
$ws.Range("I58").Value = 328.5
$ws.Range("K58").Value = 985.5
$ws.Range("M58").Value = -835.5

$ws.Range("H98").Value = 842.2759
$ws.Range("I98").Value = 1014.9474
$ws.Range("J98").Value = 514.2
$ws.Range("K98").Value = 1014.9474
$ws.Range("L98").Value = 514.2
$ws.Range("M98").Value = 483.0526
$ws.Range("N98").Value = -3510.2

$ws.Range("H122").Value = 842.2759
$ws.Range("I122").Value = 1014.9474
$ws.Range("J122").Value = 514.2
$ws.Range("K122").Value = 3044.8422
$ws.Range("L122").Value = 1542.6
$ws.Range("M122").Value = -594.8422
$ws.Range("N122").Value = -6442.6

$ws.Range("H125").Value = 1999.8334
$ws.Range("I125").Value = 0
$ws.Range("J125").Value = 1999.8334
$ws.Range("K125").Value = 0
$ws.Range("L125").Value = 17998.5006
$ws.Range("N125").Value = -22918.5006
$ws.Range("M125").ClearContents()

$ws.Range("H132").Value = 1291.7333
$ws.Range("I132").Value = 1291.7333
$ws.Range("K132").Value = 3875.199900000001
$ws.Range("M132").Value = -1345.199900000001

$ws.Range("H137").Value = 1733
$ws.Range("I137").Value = 1459.3823
$ws.Range("J137").Value = 2663.3
$ws.Range("K137").Value = 4378.1469
$ws.Range("L137").Value = 7989.900000000001
$ws.Range("M137").Value = -1828.1469
$ws.Range("N137").Value = -13089.9

$ws.Range("H138").Value = 4516.3687
$ws.Range("I138").Value = 3977.5
$ws.Range("J138").Value = 4830.7085
$ws.Range("K138").Value = 11932.5
$ws.Range("L138").Value = 14492.1255
$ws.Range("M138").Value = -6792.5
$ws.Range("N138").Value = -24772.1255

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 2168.9092
$ws.Range("I2").Value = 2052
$ws.Range("J2").Value = 2373.5
$ws.Range("K2").Value = 2052
$ws.Range("L2").Value = 2373.5
$ws.Range("M2").Value = -1939
$ws.Range("N2").Value = -2599.5

$ws.Range("H45").Value = 2945.0715
$ws.Range("I45").Value = 2627.6667
$ws.Range("J45").Value = 3516.4
$ws.Range("K45").Value = 2627.6667
$ws.Range("L45").Value = 3516.4
$ws.Range("M45").Value = -2250.6667
$ws.Range("N45").Value = -4270.4

$ws.Range("H116").Value = 2168.9092
$ws.Range("I116").Value = 2052
$ws.Range("J116").Value = 2373.5
$ws.Range("K116").Value = 2052
$ws.Range("L116").Value = 2373.5
$ws.Range("M116").Value = 242
$ws.Range("N116").Value = -6961.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 2168.9092
$ws.Range("I3").Value = 2052
$ws.Range("J3").Value = 2373.5
$ws.Range("K3").Value = 2052
$ws.Range("L3").Value = 2373.5
$ws.Range("M3").Value = -1938
$ws.Range("N3").Value = -2601.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 2999
$ws.Range("I58").Value = 2999
$ws.Range("K58").Value = 2999
$ws.Range("M58").Value = -2796

$ws.Range("H86").Value = 20911.955
$ws.Range("I86").Value = 10028.083
$ws.Range("J86").Value = 33972.6
$ws.Range("K86").Value = 10028.083
$ws.Range("L86").Value = 33972.6
$ws.Range("M86").Value = -8905.083000000001
$ws.Range("N86").Value = -36218.6

$ws.Range("H89").Value = 20911.955
$ws.Range("I89").Value = 10028.083
$ws.Range("J89").Value = 33972.6
$ws.Range("K89").Value = 50140.415
$ws.Range("L89").Value = 169863
$ws.Range("M89").Value = -44524.415
$ws.Range("N89").Value = -181095

$ws.Range("H99").Value = 9213
$ws.Range("I99").Value = 8873
$ws.Range("K99").Value = 8873
$ws.Range("M99").Value = -7375

$ws.Range("H122").Value = 2605.75
$ws.Range("J122").Value = 2605.75
$ws.Range("L122").Value = 7817.25
$ws.Range("N122").Value = -12717.25

$ws.Range("H126").Value = 9213
$ws.Range("I126").Value = 8873
$ws.Range("K126").Value = 26619
$ws.Range("M126").Value = -24149

$ws.Range("H132").Value = 4141.294
$ws.Range("I132").Value = 4282.4287
$ws.Range("K132").Value = 12847.2861
$ws.Range("M132").Value = -10317.2861

$ws.Range("H136").Value = 2999
$ws.Range("I136").Value = 2999
$ws.Range("K136").Value = 8997
$ws.Range("M136").Value = -6447

$ws.Range("H141").Value = 60710.6
$ws.Range("J141").Value = 60710.6
$ws.Range("L141").Value = 60710.6
$ws.Range("N141").Value = -71070.60000000001

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H121").Value = 3379.4
$ws.Range("I121").Value = 120
$ws.Range("J121").Value = 4194.25
$ws.Range("K121").Value = 360
$ws.Range("L121").Value = 12582.75
$ws.Range("M121").Value = 950
$ws.Range("N121").Value = -15202.75

$ws.Range("H131").Value = 1943.0667
$ws.Range("J131").Value = 1982.1428
$ws.Range("L131").Value = 5946.428400000001
$ws.Range("N131").Value = -16026.4284

$ws.Range("H132").Value = 4078.2
$ws.Range("J132").Value = 4317.4
$ws.Range("L132").Value = 38856.6
$ws.Range("N132").Value = -43916.6

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 773.7083
$ws.Range("I97").Value = 639.9474
$ws.Range("K97").Value = 639.9474
$ws.Range("M97").Value = -143.9474

$ws.Range("H107").Value = 530.625
$ws.Range("I107").Value = 540.8333
$ws.Range("J107").Value = 500
$ws.Range("K107").Value = 540.8333
$ws.Range("L107").Value = 500
$ws.Range("M107").Value = 1379.1667
$ws.Range("N107").Value = -4340

$ws.Range("H122").Value = 4216.5
$ws.Range("I122").Value = 3459.8
$ws.Range("K122").Value = 10379.4
$ws.Range("M122").Value = -7929.400000000001

$ws.Range("H126").Value = 5482.4287
$ws.Range("I126").Value = 4376.2
$ws.Range("J126").Value = 8248
$ws.Range("K126").Value = 13128.6
$ws.Range("L126").Value = 24744
$ws.Range("M126").Value = -10658.6
$ws.Range("N126").Value = -29684

$ws.Range("H132").Value = 4639.8237
$ws.Range("I132").Value = 4643.727
$ws.Range("J132").Value = 4632.6665
$ws.Range("K132").Value = 13931.181
$ws.Range("L132").Value = 13897.9995
$ws.Range("M132").Value = -11401.181
$ws.Range("N132").Value = -18957.9995

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 0
$ws.Range("I7").Value = 0
$ws.Range("K7").Value = 0
$ws.Range("M7").ClearContents()

$ws.Range("H22").Value = 5967.8
$ws.Range("J22").Value = 6316.1665
$ws.Range("L22").Value = 6316.1665
$ws.Range("N22").Value = -6906.1665

$ws.Range("H27").Value = 5967.8
$ws.Range("J27").Value = 6316.1665
$ws.Range("L27").Value = 6316.1665
$ws.Range("N27").Value = -6530.1665

$ws.Range("H61").Value = 3915
$ws.Range("I61").Value = 3956.1667
$ws.Range("J61").Value = 3832.6667
$ws.Range("K61").Value = 3956.1667
$ws.Range("L61").Value = 3832.6667
$ws.Range("M61").Value = -3754.1667
$ws.Range("N61").Value = -4236.6667

$ws.Range("H113").Value = 3915
$ws.Range("I113").Value = 3956.1667
$ws.Range("J113").Value = 3832.6667
$ws.Range("K113").Value = 3956.1667
$ws.Range("L113").Value = 3832.6667
$ws.Range("M113").Value = -1786.1667
$ws.Range("N113").Value = -8172.6667

$ws.Range("H126").Value = 0
$ws.Range("I126").Value = 0
$ws.Range("K126").Value = 0
$ws.Range("M126").ClearContents()

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H54").Value = 16337.5
$ws.Range("J54").Value = 16337.5
$ws.Range("L54").Value = 16337.5
$ws.Range("N54").Value = -17377.5

$ws.Range("H126").Value = 2119.4167
$ws.Range("I126").Value = 2175.7273
$ws.Range("K126").Value = 6527.1819
$ws.Range("M126").Value = -4057.1819

Write-Host "Updated profit figures across ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR"
